{"js": "// Each cell of the single practice-problems table holds one arithmetic\n// expression (e.g. \"5-4=1\"). The commit swaps every expression for a new\n// one while leaving the table layout, fonts and paragraph formatting\n// untouched. `replacements` lists [oldText, newText] pairs in the same\n// left-to-right / top-to-bottom order the cells appear in the table.\nconst replacements = [[\"5-4=1\", \"51+4=55\"], [\"40+10=50\", \"23+57=80\"], [\"6+78=84\", \"62+10=72\"], [\"99-18=81\", \"85-12=73\"], [\"19-0=19\", \"60-41=19\"], [\"38+61=99\", \"4+70=74\"], [\"90-75=15\", \"82-20=62\"], [\"4+20=24\", \"81-11=70\"], [\"59-55=4\", \"46-21=25\"], [\"74-69=5\", \"86-63=23\"], [\"16+39=55\", \"62-20=42\"], [\"1+29=30\", \"97-28=69\"], [\"23+12=35\", \"17+75=92\"], [\"36-6=30\", \"8+8=16\"], [\"20+24=44\", \"67+13=80\"], [\"88-18=70\", \"63-56=7\"], [\"45+6=51\", \"76+16=92\"], [\"20+27=47\", \"5+88=93\"], [\"68+2=70\", \"57+17=74\"], [\"45-33=12\", \"38-20=18\"], [\"76-60=16\", \"22+71=93\"], [\"27+70=97\", \"97-59=38\"], [\"62-32=30\", \"73-1=72\"], [\"26+64=90\", \"7+55=62\"], [\"33+51=84\", \"7-1=6\"], [\"48+31=79\", \"56+25=81\"], [\"17-7=10\", \"35+42=77\"], [\"47+46=93\", \"90+3=93\"], [\"92-74=18\", \"12+67=79\"], [\"90-12=78\", \"73+17=90\"], [\"78-49=29\", \"46-40=6\"], [\"14+38=52\", \"3+84=87\"], [\"79-26=53\", \"87-72=15\"], [\"68-28=40\", \"28+44=72\"], [\"73-8=65\", \"75-2=73\"], [\"81-34=47\", \"16+21=37\"], [\"75+5=80\", \"54+26=80\"], [\"95-75=20\", \"83-69=14\"], [\"50+45=95\", \"96-36=60\"], [\"80-50=30\", \"85-34=51\"], [\"37-28=9\", \"72-69=3\"], [\"31+32=63\", \"62-57=5\"], [\"85-61=24\", \"98-44=54\"], [\"75-1=74\", \"23-4=19\"], [\"73+1=74\", \"85-85=0\"], [\"54+41=95\", \"67+0=67\"], [\"56+20=76\", \"24+12=36\"], [\"12+60=72\", \"95+4=99\"], [\"94-76=18\", \"28+61=89\"], [\"35+15=50\", \"60+27=87\"], [\"29-13=16\", \"40+14=54\"], [\"47+42=89\", \"70-37=33\"], [\"76-65=11\", \"53-42=11\"], [\"19+77=96\", \"8+53=61\"], [\"38+16=54\", \"20+17=37\"], [\"13+82=95\", \"13+56=69\"], [\"7+60=67\", \"26+12=38\"], [\"37+49=86\", \"1+21=22\"], [\"60-34=26\", \"71-22=49\"], [\"47+8=55\", \"56-37=19\"], [\"86-71=15\", \"31+55=86\"], [\"52-24=28\", \"47+24=71\"], [\"65+30=95\", \"58+1=59\"], [\"33+64=97\", \"25+36=61\"], [\"77-13=64\", \"80+16=96\"], [\"82-6=76\", \"12+43=55\"], [\"14+61=75\", \"75-66=9\"], [\"47+43=90\", \"31+0=31\"], [\"48+2=50\", \"53+0=53\"], [\"38+33=71\", \"6+93=99\"], [\"50-41=9\", \"29+51=80\"], [\"51-0=51\", \"58+21=79\"], [\"48+13=61\", \"15+1=16\"], [\"25-21=4\", \"46-0=46\"], [\"53+40=93\", \"2+69=71\"], [\"87-19=68\", \"25+56=81\"], [\"19+15=34\", \"15+54=69\"], [\"10-7=3\", \"60-30=30\"], [\"30-7=23\", \"66+2=68\"], [\"48+25=73\", \"1+3=4\"], [\"81-67=14\", \"92+7=99\"], [\"58-37=21\", \"99-24=75\"], [\"48+49=97\", \"28-19=9\"], [\"61-1=60\", \"63+17=80\"], [\"25+59=84\", \"65-38=27\"], [\"90-16=74\", \"46-16=30\"], [\"20+50=70\", \"13+55=68\"], [\"34+34=68\", \"33+8=41\"], [\"10+40=50\", \"52-14=38\"], [\"8+29=37\", \"76-7=69\"], [\"91-36=55\", \"40+18=58\"], [\"48+13=61\", \"44-15=29\"], [\"77+9=86\", \"26+30=56\"], [\"99-44=55\", \"30+69=99\"], [\"54+25=79\", \"38+29=67\"], [\"76-34=42\", \"98-92=6\"], [\"99-51=48\", \"45+49=94\"], [\"46-11=35\", \"39-23=16\"], [\"86-65=21\", \"38+5=43\"], [\"20-2=18\", \"51+39=90\"]];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\nawait context.sync();\nif (table.isNullObject) {\n  throw new Error(\"Expected exactly one table in the document body.\");\n}\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\n// Walk the cells in document order, matching each one's current text\n// against the next unconsumed replacement pair (falls back to a text\n// lookup if the order ever doesn't line up, so this is resilient either\n// way).\nlet cursor = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const current = cell.body.text.replace(/[\\r\\v\\f]+$/g, \"\").trim();\n    let newText = null;\n    if (cursor < replacements.length && replacements[cursor][0] === current) {\n      newText = replacements[cursor][1];\n      cursor++;\n    } else {\n      const hit = replacements.find(([oldText]) => oldText === current);\n      if (hit) {\n        newText = hit[1];\n      }\n    }\n    if (newText !== null) {\n      cell.value = newText;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document body contains a single 20x5 table where every cell holds\n# one arithmetic expression (e.g. \"5-4=1\"). The commit swaps every\n# expression for a new one while leaving the table layout, fonts and\n# paragraph formatting untouched. Addressing cells by (row, column) keeps\n# this correct even though a couple of the old expressions repeat verbatim\n# (e.g. \"48+13=61\" appears twice but maps to two different new values),\n# which a single document-wide Find/Replace could not distinguish.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New expression for every cell, in row-major order (row 1..20, col 1..5),\n# matching the table's reading order top-to-bottom / left-to-right.\n$values = @(\n    @(\"51+4=55\", \"23+57=80\", \"62+10=72\", \"85-12=73\", \"60-41=19\"),\n    @(\"4+70=74\", \"82-20=62\", \"81-11=70\", \"46-21=25\", \"86-63=23\"),\n    @(\"62-20=42\", \"97-28=69\", \"17+75=92\", \"8+8=16\", \"67+13=80\"),\n    @(\"63-56=7\", \"76+16=92\", \"5+88=93\", \"57+17=74\", \"38-20=18\"),\n    @(\"22+71=93\", \"97-59=38\", \"73-1=72\", \"7+55=62\", \"7-1=6\"),\n    @(\"56+25=81\", \"35+42=77\", \"90+3=93\", \"12+67=79\", \"73+17=90\"),\n    @(\"46-40=6\", \"3+84=87\", \"87-72=15\", \"28+44=72\", \"75-2=73\"),\n    @(\"16+21=37\", \"54+26=80\", \"83-69=14\", \"96-36=60\", \"85-34=51\"),\n    @(\"72-69=3\", \"62-57=5\", \"98-44=54\", \"23-4=19\", \"85-85=0\"),\n    @(\"67+0=67\", \"24+12=36\", \"95+4=99\", \"28+61=89\", \"60+27=87\"),\n    @(\"40+14=54\", \"70-37=33\", \"53-42=11\", \"8+53=61\", \"20+17=37\"),\n    @(\"13+56=69\", \"26+12=38\", \"1+21=22\", \"71-22=49\", \"56-37=19\"),\n    @(\"31+55=86\", \"47+24=71\", \"58+1=59\", \"25+36=61\", \"80+16=96\"),\n    @(\"12+43=55\", \"75-66=9\", \"31+0=31\", \"53+0=53\", \"6+93=99\"),\n    @(\"29+51=80\", \"58+21=79\", \"15+1=16\", \"46-0=46\", \"2+69=71\"),\n    @(\"25+56=81\", \"15+54=69\", \"60-30=30\", \"66+2=68\", \"1+3=4\"),\n    @(\"92+7=99\", \"99-24=75\", \"28-19=9\", \"63+17=80\", \"65-38=27\"),\n    @(\"46-16=30\", \"13+55=68\", \"33+8=41\", \"52-14=38\", \"76-7=69\"),\n    @(\"40+18=58\", \"44-15=29\", \"26+30=56\", \"30+69=99\", \"38+29=67\"),\n    @(\"98-92=6\", \"45+49=94\", \"39-23=16\", \"38+5=43\", \"51+39=90\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $values[$r - 1][$c - 1]\n    }\n}\n"}
